# Update reports - 2026-01-27 10:10
#
# A fresh scan ran (2026-01-27 10:10:11) and found no new publication IDs,
# so only a "Scans" log row is appended. The same re-scan re-affirmed the
# four existing "Order Papers" rows, so they are duplicated as new rows.
# Both tables are then resized to cover their newly added rows.

$wb = $excel.ActiveWorkbook

# Helper: write a literal piece of text into a cell even when the text
# looks like a date/time (e.g. "2026-01-27", "10:10:11") so Excel doesn't
# silently convert it to a date/time serial number. A leading apostrophe
# forces text entry (exactly like typing '2026-01-27 into the formula
# bar); re-applying the "Normal" style afterwards drops the resulting
# quote-prefix formatting so the cell is left with the default style.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "Scans": append row 3 -> scan at 10:10:11 found no new IDs
# ---------------------------------------------------------------------
$wsScans = $wb.Worksheets.Item("Scans")

Set-TextValue $wsScans.Range("A3") "2026-01-27"
Set-TextValue $wsScans.Range("B3") "10:10:11"

$scansTable = $wsScans.ListObjects.Item("Scans")
$scansTable.Resize($wsScans.Range("A1:C3"))

# ---------------------------------------------------------------------
# Sheet "Order Papers": duplicate rows 2-5 into new rows 6-9
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Order Papers")

# Row 6 = copy of row 2
Set-TextValue $wsOrders.Range("A6") "2026-01-22"
$wsOrders.Range("B6").Value = "Health and Social Care"
$wsOrders.Range("C6").Value = "5th Report: First 1000 Days: a renewed focus"
$wsOrders.Range("D6").Value = "HC 802"
Set-TextValue $wsOrders.Range("E6") "2026-01-22"
Set-TextValue $wsOrders.Range("F6") "00:01:00"

# Row 7 = copy of row 3
Set-TextValue $wsOrders.Range("A7") "2026-01-22"
$wsOrders.Range("B7").Value = "International Development"
$wsOrders.Range("C7").Value = "7th Special Report: Empowering Development: Energy Access for Communities: Government Response"
$wsOrders.Range("D7").Value = "HC 1626"
Set-TextValue $wsOrders.Range("E7") "2026-01-22"
Set-TextValue $wsOrders.Range("F7") "00:01:00"

# Row 8 = copy of row 4
Set-TextValue $wsOrders.Range("A8") "2026-01-22"
$wsOrders.Range("B8").Value = "Treasury"
$wsOrders.Range("C8").Value = "6th Special Report: Taxation of gambling: Government Response"
$wsOrders.Range("D8").Value = "HC 1625"
Set-TextValue $wsOrders.Range("E8") "2026-01-22"
Set-TextValue $wsOrders.Range("F8") "00:01:00"

# Row 9 = copy of row 5
Set-TextValue $wsOrders.Range("A9") "2026-01-22"
$wsOrders.Range("B9").Value = "Home Affairs"
$wsOrders.Range("C9").Value = "3rd Special Report: The Home Office’s management of asylum accommodation: Government Response"
$wsOrders.Range("D9").Value = "HC 1642"
Set-TextValue $wsOrders.Range("E9") "2026-01-22"
Set-TextValue $wsOrders.Range("F9") "10:00:00"

$ordersTable = $wsOrders.ListObjects.Item("Order_Papers")
$ordersTable.Resize($wsOrders.Range("A1:H9"))
